$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the header cell from "keyword" to "whack_word"
$ws.Range("A1").Value = "whack_word"

# The list had two accidental duplicate keywords ("equal opportunity" on
# row 12 and "be a" on row 8). Remove the duplicate rows, deleting from
# the bottom up so row numbers of earlier rows are unaffected.
$ws.Rows(12).Delete()
$ws.Rows(8).Delete()

# Reflect that the user finished up by selecting the whole column and
# scrolling down toward the bottom of the (now shorter) list before saving.
$ws.Columns("A:A").Select()
$excel.ActiveWindow.ScrollRow = 34
$excel.ActiveWindow.ScrollColumn = 1
